# Commit: "Provide service to upload spreadsheets; rename ContactRowHandler to ContactImport"
#
# The only deliberate, content-level change captured by the diff for this
# workbook is renaming the single worksheet from "ContactRowHandler" to
# "ContactImport" (the rest of the hunks - fileVersion/rupBuild, absPath,
# window geometry, calcPr/concurrentCalc, the mx:ArchID extLst, the
# x14ac:dyDescent/default row height/column bestFit width tweaks and the
# East-Asian theme font substitutions - are all side effects of the file
# having been re-saved from a different Excel build/platform, not
# something a user did through the object model).
#
# While editing, the sheet's selection also ended up on cell E2, so we
# replicate that too since it is directly expressed in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "ContactImport"
$ws.Range("E2").Select()
